# Append 4 new incident-log rows (171-174) to Sheet1, extending the
# existing A:G table of Fecha/Hora/WC47.../SPL entries.
# Columns A (Fecha) hold date-like text ("2024-05-14"); Excel would
# normally auto-convert such strings to date serials on plain .Value
# assignment, so those cells are pre-formatted as Text ("@") before the
# value is written, matching how the source file stores them as literal
# strings rather than numeric dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 171; A = "2024-05-14"; B = "18:56:57"; C = "-";                  D = "Tornillo atascado en tolva"; E = "-";              F = "-"; G = "-" },
    @{ Row = 172; A = "2024-05-14"; B = "18:57:01"; C = "-";                  D = "Fallo etiqueta";              E = "-";              F = "-"; G = "-" },
    @{ Row = 173; A = "2024-05-14"; B = "19:23:33"; C = "Fallo en elevador";  D = "-";                           E = "-";              F = "-"; G = "-" },
    @{ Row = 174; A = "2024-05-14"; B = "19:31:47"; C = "-";                  D = "-";                           E = "Atasco tuerca";  F = "-"; G = "-" }
)

foreach ($r in $rows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.A

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
